# Auto-generated Excel COM-interop edit script
# Applies the cell-value changes described in the commit diff
# (profit/price recalculations across several Leve sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17 (Leve Item ID 38956)
$ws.Range("H17").Value = 202.5
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 202.5
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 607.5
$ws.Range("M17").Value = $null
$ws.Range("N17").Value = -943.5

# Row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 1689.6666
$ws.Range("I40").Value = 1618.7142
$ws.Range("J40").Value = 1789
$ws.Range("K40").Value = 1618.7142
$ws.Range("L40").Value = 1789
$ws.Range("M40").Value = -1443.7142
$ws.Range("N40").Value = -2139

# Row 126 (Leve Item ID 34391)
$ws.Range("H126").Value = 39712.727
$ws.Range("J126").Value = 39712.727
$ws.Range("L126").Value = 39712.727
$ws.Range("N126").Value = -49592.727

# Row 135 (Leve Item ID 44047)
$ws.Range("H135").Value = 4759.3184
$ws.Range("I135").Value = 3578.111
$ws.Range("J135").Value = 5577.077
$ws.Range("K135").Value = 32202.999
$ws.Range("L135").Value = 50193.693
$ws.Range("M135").Value = -29667.999
$ws.Range("N135").Value = -55263.693

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 1279.1666
$ws.Range("I137").Value = 975
$ws.Range("J137").Value = 1431.25
$ws.Range("K137").Value = 2925
$ws.Range("L137").Value = 4293.75
$ws.Range("M137").Value = -375
$ws.Range("N137").Value = -9393.75

$ws = $wb.Worksheets.Item("ARM")
# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 2277.375
$ws.Range("I45").Value = 2124.6785
$ws.Range("J45").Value = 2633.6667
$ws.Range("K45").Value = 2124.6785
$ws.Range("L45").Value = 2633.6667
$ws.Range("M45").Value = -1747.6785
$ws.Range("N45").Value = -3387.6667

# Row 141 (Leve Item ID 42483)
$ws.Range("H141").Value = 65589.31
$ws.Range("J141").Value = 65589.31
$ws.Range("L141").Value = 65589.31
$ws.Range("N141").Value = -75949.31

$ws = $wb.Worksheets.Item("BSM")
# Row 132 (Leve Item ID 41855)
$ws.Range("H132").Value = 8500
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = $null

$ws = $wb.Worksheets.Item("CRP")
# Row 7 (Leve Item ID 5361)
$ws.Range("H7").Value = 48.090908
$ws.Range("I7").Value = 38.375
$ws.Range("J7").Value = 74
$ws.Range("K7").Value = 38.375
$ws.Range("L7").Value = 74
$ws.Range("M7").Value = 74.625
$ws.Range("N7").Value = -300

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (Leve Item ID 43974)
$ws.Range("H5").Value = 948.6875
$ws.Range("I5").Value = 1009.2143
$ws.Range("J5").Value = 525
$ws.Range("K5").Value = 3027.6429
$ws.Range("L5").Value = 1575
$ws.Range("M5").Value = -2915.6429
$ws.Range("N5").Value = -1799

# Row 86 (Leve Item ID 12892)
$ws.Range("H86").Value = 1430
$ws.Range("J86").Value = 1430
$ws.Range("L86").Value = 4290
$ws.Range("N86").Value = -6662

# Row 89 (Leve Item ID 12892)
$ws.Range("H89").Value = 1430
$ws.Range("J89").Value = 1430
$ws.Range("L89").Value = 12870
$ws.Range("N89").Value = -24726

# Row 129 (Leve Item ID 36054)
$ws.Range("H129").Value = 13335301
$ws.Range("I129").Value = 2486
$ws.Range("J129").Value = 16668505
$ws.Range("K129").Value = 7458
$ws.Range("L129").Value = 50005515
$ws.Range("M129").Value = -2458
$ws.Range("N129").Value = -50015515

# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 6098547.5
$ws.Range("I131").Value = 1510
$ws.Range("J131").Value = 6494459
$ws.Range("K131").Value = 4530
$ws.Range("L131").Value = 19483377
$ws.Range("M131").Value = 510
$ws.Range("N131").Value = -19493457

# Row 135 (Leve Item ID 43974)
$ws.Range("H135").Value = 948.6875
$ws.Range("I135").Value = 1009.2143
$ws.Range("J135").Value = 525
$ws.Range("K135").Value = 9082.9287
$ws.Range("L135").Value = 4725
$ws.Range("M135").Value = -6547.9287
$ws.Range("N135").Value = -9795

$ws = $wb.Worksheets.Item("GSM")
# Row 2 (Leve Item ID 5062)
$ws.Range("H2").Value = 13.166667
$ws.Range("I2").Value = 13.8
$ws.Range("K2").Value = 13.8
$ws.Range("M2").Value = 99.2

# Row 18 (Leve Item ID 4309)
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = $null
$ws.Range("N18").Value = $null

# Row 21 (Leve Item ID 4430)
$ws.Range("H21").Value = 2001600
$ws.Range("J21").Value = 2000
$ws.Range("L21").Value = 2000
$ws.Range("N21").Value = -2346

# Row 30 (Leve Item ID 4430)
$ws.Range("H30").Value = 2001600
$ws.Range("J30").Value = 2000
$ws.Range("L30").Value = 2000
$ws.Range("N30").Value = -2210

# Row 43 (Leve Item ID 4218)
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = $null

# Row 46 (Leve Item ID 2078)
$ws.Range("H46").Value = 11031.111
$ws.Range("I46").Value = 3575
$ws.Range("J46").Value = 16996
$ws.Range("K46").Value = 3575
$ws.Range("L46").Value = 16996
$ws.Range("M46").Value = -3419
$ws.Range("N46").Value = -17308

# Row 57 (Leve Item ID 2876)
$ws.Range("H57").Value = 30555.445
$ws.Range("J57").Value = 33124.875
$ws.Range("L57").Value = 33124.875
$ws.Range("N57").Value = -34764.875

# Row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 25574944
$ws.Range("I70").Value = 40183450
$ws.Range("J70").Value = 10062.5
$ws.Range("K70").Value = 40183450
$ws.Range("L70").Value = 10062.5
$ws.Range("M70").Value = -40183180
$ws.Range("N70").Value = -10602.5

# Row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 25574944
$ws.Range("I73").Value = 40183450
$ws.Range("J73").Value = 10062.5
$ws.Range("K73").Value = 40183450
$ws.Range("L73").Value = 10062.5
$ws.Range("M73").Value = -40182514
$ws.Range("N73").Value = -11934.5

# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 3763.125
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 4017.5
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 4017.5
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -6013.5

# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 3763.125
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 4017.5
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 20087.5
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -30071.5

$ws = $wb.Worksheets.Item("LTW")
# Row 68 (Leve Item ID 12563)
$ws.Range("H68").Value = 2096.862
$ws.Range("I68").Value = 2014.2858
$ws.Range("J68").Value = 2173.9333
$ws.Range("K68").Value = 2014.2858
$ws.Range("L68").Value = 2173.9333
$ws.Range("M68").Value = -1265.2858
$ws.Range("N68").Value = -3671.9333

# Row 71 (Leve Item ID 12563)
$ws.Range("H71").Value = 2096.862
$ws.Range("I71").Value = 2014.2858
$ws.Range("J71").Value = 2173.9333
$ws.Range("K71").Value = 10071.429
$ws.Range("L71").Value = 10869.6665
$ws.Range("M71").Value = -6327.429
$ws.Range("N71").Value = -18357.6665

# Row 133 (Leve Item ID 41903)
$ws.Range("H133").Value = 92626.92
$ws.Range("J133").Value = 92626.92
$ws.Range("L133").Value = 92626.92
$ws.Range("N133").Value = -97686.92

# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 18204.666
$ws.Range("I136").Value = 22030.8
$ws.Range("K136").Value = 66092.39999999999
$ws.Range("M136").Value = -63542.39999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 15 (Leve Item ID 2670)
$ws.Range("H15").Value = 8006.8
$ws.Range("J15").Value = 8006.8
$ws.Range("L15").Value = 8006.8
$ws.Range("N15").Value = -8582.799999999999

# Row 54 (Leve Item ID 3413)
$ws.Range("H54").Value = 6199
$ws.Range("J54").Value = 6199
$ws.Range("L54").Value = 6199
$ws.Range("N54").Value = -7239

# Row 81 (Leve Item ID 12596)
$ws.Range("H81").Value = 2579.24
$ws.Range("I81").Value = 1987.1177
$ws.Range("J81").Value = 3837.5
$ws.Range("K81").Value = 3974.2354
$ws.Range("L81").Value = 7675
$ws.Range("M81").Value = -2913.2354
$ws.Range("N81").Value = -9797

# Row 84 (Leve Item ID 12596)
$ws.Range("H84").Value = 2579.24
$ws.Range("I84").Value = 1987.1177
$ws.Range("J84").Value = 3837.5
$ws.Range("K84").Value = 19871.177
$ws.Range("L84").Value = 38375
$ws.Range("M84").Value = -14567.177
$ws.Range("N84").Value = -48983

# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 21681.312
$ws.Range("I136").Value = 36399.75
$ws.Range("J136").Value = 1075.5
$ws.Range("K136").Value = 109199.25
$ws.Range("L136").Value = 3226.5
$ws.Range("M136").Value = -106649.25
$ws.Range("N136").Value = -8326.5
